# Refresh crypto price/volume data to match the latest scrape.
# Numeric-looking "Price" strings (e.g. "1.20", "0.0220") must stay text,
# so they are entered with a leading apostrophe (forces text, like typing
# into Excel) and then the style is reset to a plain sibling cell so no
# quote-prefix formatting lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("D4").Style

# Row 2
$ws.Range("D2").Value = '37.423.26'
$ws.Range("E2").Value = '  +4.47%  '

# Row 3
$ws.Range("D3").Value = '2.045.41'
$ws.Range("E3").Value = '  +3.02%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '''252.33'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  +2.86%  '

# Row 6
$ws.Range("D6").Value = '''0.648'
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = '  +1.44%  '

# Row 7
$ws.Range("D7").Value = '''66.05'
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = '  +10.52%  '

# Row 9
$ws.Range("D9").Value = '''0.403'
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = '  +11.02%  '

# Row 10
$ws.Range("D10").Value = '''59.46'
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = '  +2.40%  '

# Row 11
$ws.Range("D11").Value = '''0.0809'
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  +9.31%  '

# Row 12
$ws.Range("E12").Value = '  +0.10%  '

# Row 13
$ws.Range("E13").Value = '  -2.52%  '

# Row 14
$ws.Range("D14").Value = '''23.58'
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = '  +24.32%  '

# Row 15
$ws.Range("D15").Value = '''14.82'
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  +1.08%  '

# Row 16
$ws.Range("D16").Value = '2.343.18'

# Row 17
$ws.Range("D17").Value = '''5.71'
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = '  +7.07%  '

# Row 18
$ws.Range("D18").Value = '2.049.76'
$ws.Range("E18").Value = '  +3.46%  '

# Row 19
$ws.Range("D19").Value = '37.257.18'
$ws.Range("E19").Value = '  +4.23%  '

# Row 20
$ws.Range("D20").Value = '''73.22'
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = '  +2.35%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0902'
$ws.Range("E21").Value = '  +6.50%  '

# Row 22
$ws.Range("D22").Value = '''5.52'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  +6.13%  '

# Row 23
$ws.Range("D23").Value = '''239.13'
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = '  +2.90%  '

# Row 24
$ws.Range("E24").Value = '  -0.05%  '

# Row 25
$ws.Range("D25").Value = '''2.61'
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  +1.05%  '

# Row 26
$ws.Range("E26").Value = '  +3.94%  '

# Row 27
$ws.Range("D27").Value = '''10.07'
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = '  +4.40%  '

# Row 28
$ws.Range("D28").Value = '''162.05'
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  -1.74%  '

# Row 29
$ws.Range("D29").Value = '''20.13'
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = '  +3.95%  '

# Row 30
$ws.Range("D30").Value = '''0.131'
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = '  +33.11%  '

# Row 31
$ws.Range("E31").Value = '  +3.04%  '

# Row 32
$ws.Range("D32").Value = '''5.16'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  +3.13%  '

# Row 33
$ws.Range("D33").Value = '''1.20'
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = '  +5.35%  '

# Row 34
$ws.Range("E34").Value = '  +4.85%  '

# Row 35
$ws.Range("D35").Value = '''4.67'
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = '  +4.87%  '

# Row 36
$ws.Range("E36").Value = '  +13.19%  '

# Row 37
$ws.Range("D37").Value = '''2.38'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  -2.25%  '

# Row 38
$ws.Range("E38").Value = '  +0.02%  '

# Row 39
$ws.Range("E39").Value = '  +3.96%  '

# Row 40
$ws.Range("D40").Value = '''3.02'
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = '  +31.14%  '

# Row 41
$ws.Range("D41").Value = '''1.31'
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  +6.46%  '

# Row 42
$ws.Range("D42").Value = '''0.102'
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  +8.61%  '

# Row 43
$ws.Range("E43").Value = '  +5.64%  '

# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''17.54'
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = '  +6.33%  '

# Row 45
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = '''1.16'
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = '  +5.56%  '

# Row 46
$ws.Range("D46").Value = '''0.0220'
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  +2.76%  '

# Row 47
$ws.Range("D47").Value = '''95.58'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  +2.69%  '

# Row 48
$ws.Range("E48").Value = '  +1.20%  '

# Row 49
$ws.Range("D49").Value = '1.394.81'
$ws.Range("E49").Value = '  +2.62%  '

# Row 50
$ws.Range("D50").Value = '''2.93'
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  +1.44%  '

# Row 51
$ws.Range("D51").Value = '''47.29'
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = '  +1.57%  '
